$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.842.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.937.78'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.51'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4909'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2936'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06881'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '104.89'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07773'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.931.72'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.352'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6985'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '273.98'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.857.68'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007718'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.08'
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.198.57'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.572'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.527'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.859'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.07'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.60'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.158'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1042'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.25%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.556'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.550'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.376'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04886'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7580'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.150'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.80%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.656'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.538'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.70'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.091'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9124'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4435'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.79'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9988'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.693'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -7.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '995.32'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.50%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.09'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.16%  '
